$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 13 de Mayo de 2020 a las 18:05"

# Estados Unidos (row 4) - refreshed totals
$ws.Range("B4").Value = 1413190
$ws.Range("C4").Value = 4554
$ws.Range("D4").Value = 299008
$ws.Range("E4").Value = 1030469
$ws.Range("F4").Value = 16473
$ws.Range("G4").Value = 288
$ws.Range("H4").Value = 83713

# Canada (row 17) - refreshed totals
$ws.Range("B17").Value = 71486
$ws.Range("C17").Value = 329
$ws.Range("D17").Value = 34496
$ws.Range("E17").Value = 31781
$ws.Range("F17").Value = 502
$ws.Range("G17").Value = 40
$ws.Range("H17").Value = 5209

# Chile overtakes Pakistan: Chile now sits at row 22, Pakistan moves to row 23
$ws.Range("A22").Value = "Chile"
$ws.Range("B22").Value = 34381
$ws.Range("C22").Value = 2660
$ws.Range("D22").Value = 14865
$ws.Range("E22").Value = 19170
$ws.Range("F22").Value = 494
$ws.Range("G22").Value = 11
$ws.Range("H22").Value = 346

$ws.Range("A23").Value = "Pakistan"
$ws.Range("B23").Value = 34336
$ws.Range("C23").Value = 1662
$ws.Range("D23").Value = 8812
$ws.Range("E23").Value = 24787
$ws.Range("F23").Value = 111
$ws.Range("G23").Value = 13
$ws.Range("H23").Value = 737

# Polonia (row 34) - refreshed totals
$ws.Range("B34").Value = 17204
$ws.Range("C34").Value = 283
$ws.Range("D34").Value = 6410
$ws.Range("E34").Value = 9933
$ws.Range("F34").Value = 160
$ws.Range("G34").Value = 22
$ws.Range("H34").Value = 861

# Republica Dominicana overtakes Kuwait & Corea del Sur:
# Republica Dominicana now sits at row 44, Kuwait moves to row 45, Corea del Sur moves to row 46
$ws.Range("A44").Value = "Republica Dominicana"
$ws.Range("B44").Value = 11196
$ws.Range("C44").Value = 296
$ws.Range("D44").Value = 3221
$ws.Range("E44").Value = 7566
$ws.Range("F44").Value = 131
$ws.Range("G44").Value = 7
$ws.Range("H44").Value = 409

$ws.Range("A45").Value = "Kuwait"
$ws.Range("B45").Value = 11028
$ws.Range("C45").Value = 751
$ws.Range("D45").Value = 3263
$ws.Range("E45").Value = 7683
$ws.Range("F45").Value = 169
$ws.Range("G45").Value = 7
$ws.Range("H45").Value = 82

$ws.Range("A46").Value = "Corea del Sur"
$ws.Range("B46").Value = 10962
$ws.Range("C46").Value = 26
$ws.Range("D46").Value = 9695
$ws.Range("E46").Value = 1008
$ws.Range("F46").Value = 55
$ws.Range("G46").Value = 1
$ws.Range("H46").Value = 259

# Chequia (row 51) - refreshed totals
$ws.Range("B51").Value = 8240
$ws.Range("C51").Value = 42
$ws.Range("D51").Value = 5043
$ws.Range("E51").Value = 2909
$ws.Range("F51").Value = 43
$ws.Range("G51").Value = 5
$ws.Range("H51").Value = 288

# Luxemburgo (row 66) - refreshed totals
$ws.Range("B66").Value = 3904
$ws.Range("C66").Value = 10
$ws.Range("D66").Value = 3629
$ws.Range("E66").Value = 172
$ws.Range("F66").Value = 22
$ws.Range("G66").Value = 1
$ws.Range("H66").Value = 103

# Principado de Andorra (row 113) - refreshed totals
$ws.Range("B113").Value = 760
$ws.Range("C113").Value = 2
$ws.Range("D113").Value = 576
$ws.Range("E113").Value = 135
$ws.Range("F113").Value = 3
$ws.Range("G113").Value = 1
$ws.Range("H113").Value = 49
